# Add a new worksheet "22097" after the existing "GCHSB" sheet and fill it
# with RollCode/RollNumber data (mirrors the structure of the GCHSB sheet).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "22097"

# Header row - reuse the exact same text/format as the GCHSB header so the
# shared string table and style table are reused rather than duplicated.
$newSheet.Range("A1").Value = "RollCode"
$newSheet.Range("B1").Value = "RollNumber"
$h1 = $ws1.Range("A1:B1").Copy()
$h2 = $newSheet.Range("A1:B1").PasteSpecial(-4122)

# Data rows: column A holds the numeric roll code, column B holds the
# zero-padded roll number stored as text.
for ($i = 2; $i -le 30; $i++) {
    $newSheet.Cells.Item($i, 1).Value = 22097
    $newSheet.Range("B" + $i).NumberFormat = "@"
    $newSheet.Cells.Item($i, 2).Value = "{0:D4}" -f ($i - 1)
}

# Approximate the original "best fit" column widths.
$newSheet.Columns.Item(1).ColumnWidth = 8.17
$newSheet.Columns.Item(2).ColumnWidth = 10.83

# Page setup to match the other sheet (A4, portrait).
$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

# Restore/settle the selection on the original sheet (no longer the active
# tab) and leave a fresh selection + the active tab on the new sheet.
$s1 = $ws1.Range("A1:B1").Select()
$s2 = $newSheet.Range("G8").Select()

Write-Host "Added sheet 22097"
